$wb = $excel.ActiveWorkbook
$wsJugadores = $wb.Worksheets.Item("Jugadores")
$wsPartidos  = $wb.Worksheets.Item("Partidos")

# --- Jugadores: fill in missing "# camiseta" (column D) numbers ---
$wsJugadores.Range("D10").Value = 6
$wsJugadores.Range("D15").Value = 14
$wsJugadores.Range("D23").Value = 3

# --- Partidos: append the results for the 2025-11-08 (serial 45969) match ---
$matchDate = 45969

$newRows = @(
    @{ Jugador = "Sombra";                   Equipo = "Amarillo"; Posicion = "Arquero";        Goles = 0; Autogoles = 0; Arquero = $true;  GolesRecibidos = 2; Amarillas = 0; Rojas = 0; Asistencias = 0; Penales = 0 },
    @{ Jugador = "Fabian Caicedo";            Equipo = "Azul";     Posicion = "Arquero";        Goles = 0; Autogoles = 0; Arquero = $true;  GolesRecibidos = 5; Amarillas = 0; Rojas = 0; Asistencias = 0; Penales = 0 },
    @{ Jugador = "Juan David Espinal";        Equipo = "Amarillo"; Posicion = "Mediocampista";  Goles = 1; Autogoles = 0; Arquero = $false; GolesRecibidos = 0; Amarillas = 1; Rojas = 0; Asistencias = 0; Penales = 0 },
    @{ Jugador = "Cesar Augusto Estrada";     Equipo = "Amarillo"; Posicion = "Delantero";      Goles = 1; Autogoles = 0; Arquero = $false; GolesRecibidos = 0; Amarillas = 0; Rojas = 0; Asistencias = 0; Penales = 0 },
    @{ Jugador = "Alexander Uribe";           Equipo = "Azul";     Posicion = "Mediocampista";  Goles = 2; Autogoles = 0; Arquero = $false; GolesRecibidos = 0; Amarillas = 0; Rojas = 0; Asistencias = 1; Penales = 0 },
    @{ Jugador = "Carlos Fernando Valencia";  Equipo = "Azul";     Posicion = "Delantero";      Goles = 0; Autogoles = 0; Arquero = $false; GolesRecibidos = 0; Amarillas = 0; Rojas = 0; Asistencias = 1; Penales = 0 },
    @{ Jugador = "David Fernando Velasco";    Equipo = "Azul";     Posicion = "Delantero";      Goles = 2; Autogoles = 0; Arquero = $false; GolesRecibidos = 0; Amarillas = 0; Rojas = 0; Asistencias = 0; Penales = 0 },
    @{ Jugador = "Francisco Javier Duran";    Equipo = "Azul";     Posicion = "Defensa";        Goles = 1; Autogoles = 0; Arquero = $false; GolesRecibidos = 0; Amarillas = 0; Rojas = 0; Asistencias = 0; Penales = 0 },
    @{ Jugador = "Bryan Andres Burgos";       Equipo = "Azul";     Posicion = "Mediocampista";  Goles = 0; Autogoles = 0; Arquero = $false; GolesRecibidos = 0; Amarillas = 0; Rojas = 0; Asistencias = 1; Penales = 0 }
)

$startRow = 517
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $wsPartidos.Cells.Item($r, 1).Value  = $matchDate
    $wsPartidos.Cells.Item($r, 2).Value  = $row.Jugador
    $wsPartidos.Cells.Item($r, 3).Value  = $row.Equipo
    $wsPartidos.Cells.Item($r, 4).Value  = $row.Posicion
    $wsPartidos.Cells.Item($r, 5).Value  = $row.Goles
    $wsPartidos.Cells.Item($r, 6).Value  = $row.Autogoles
    $wsPartidos.Cells.Item($r, 7).Value  = $row.Arquero
    $wsPartidos.Cells.Item($r, 8).Value  = $row.GolesRecibidos
    $wsPartidos.Cells.Item($r, 9).Value  = $row.Amarillas
    $wsPartidos.Cells.Item($r, 10).Value = $row.Rojas
    $wsPartidos.Cells.Item($r, 11).Value = $row.Asistencias
    $wsPartidos.Cells.Item($r, 12).Value = $row.Penales
}

$lastRow = $startRow + $newRows.Count - 1

# --- Restore view/selection state to match the saved workbook ---
$wsJugadores.Range("F18").Select()
$wsPartidos.Activate()
$selRow = $lastRow + 1
$wsPartidos.Cells.Item($selRow, 1).Select()
